# Insert a new daily-log entry for 2026/01/21 (time 14) just above the
# existing row 685 ("2026/12/29"), pushing that row and every row below it
# down by one. This mirrors the diff: dimension grows from D726 to D727,
# and rows 685-726 become 686-727 with their original contents untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 685 (and everything after it) down one row, leaving a blank
# row 685 for the new record.
$ws.Rows.Item(685).Insert()

# Column A holds the date as plain text (matching the rest of the sheet,
# which stores "YYYY/MM/DD" as literal strings, not real dates). The
# leading apostrophe forces text entry so it isn't auto-converted into a
# date serial value; ClearFormats() then drops the resulting quote-prefix
# style so the cell ends up with no explicit style, same as its siblings.
$ws.Cells.Item(685, 1).Value = "'2026/01/21"
$ws.Cells.Item(685, 1).ClearFormats()

$ws.Cells.Item(685, 2).Value = "水"
$ws.Cells.Item(685, 3).Value = 14
$ws.Cells.Item(685, 4).Value = 201
